# The first exhibition entry ("丽水·CCAC动漫游戏嘉年华") was removed from the
# "展览" (Exhibitions) and "全部类型" (All Types) sheets. The remaining three
# events shift up by one data row (their B:I details move from old rows
# 3/4/5 into rows 2/3/4), while column A's sequential index numbers
# (1, 2, 3) stay exactly where they already were. The old trailing row
# (row 5) is then deleted entirely, shrinking the used range from A1:I5
# down to A1:I4.

$sheetNames = @("展览", "全部类型")

# New contents for rows 2-4, columns B-I (column A is left untouched).
$rowsData = @(
    @{
        B = "2024-05-03"
        C = "丽水·首届TCT国风动漫游戏嘉年华（取消）"
        D = "括苍路493号油泵厂山顶通用设备厂区块3号楼 中国国际摄影节展览馆"
        E = "2024.05.03 10:00-05.04 17:00"
        F = 167
        G = "不可售"
        H = "https://show.bilibili.com/platform/detail.html?id=84156"
        I = "//i2.hdslb.com/bfs/openplatform/202404/rpRBCHaU1712892375435.jpeg"
    },
    @{
        B = "2024-05-18"
        C = "丽水·第三届HP国风动漫游戏嘉年华"
        D = "好溪路与望城路交汇西北侧地块 丽水市水上运动中心"
        E = "2024.05.18 09:00-05.18 17:00"
        F = 180
        G = 68
        H = "https://show.bilibili.com/platform/detail.html?id=82901"
        I = "//i1.hdslb.com/bfs/openplatform/202403/sl5TubQI1710410535537.jpeg"
    },
    @{
        B = "2024-06-01"
        C = "丽水·动漫游戏展"
        D = "中东路848号(解放街交汇) 飞达国际大酒店"
        E = "2024.06.01 10:00-06.01 17:00"
        F = 131
        G = 45
        H = "https://show.bilibili.com/platform/detail.html?id=84450"
        I = "//i2.hdslb.com/bfs/openplatform/202404/tdhb9QSW1713333412467.jpeg"
    }
)

$wb = $excel.ActiveWorkbook

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    for ($i = 0; $i -lt $rowsData.Length; $i++) {
        $row = $i + 2
        $data = $rowsData[$i]

        # Column B holds a date-shaped string ("2024-05-03"). Force the
        # cell to Text format first so Excel does not silently convert it
        # into a real date serial number, then restore the default
        # (unstyled) cell style so no numFmt/style index leaks into the
        # saved file.
        $bCell = $ws.Cells.Item($row, 2)
        $bCell.NumberFormat = "@"
        $bCell.Value = $data.B
        $bCell.Style = "Normal"

        $ws.Cells.Item($row, 3).Value = $data.C
        $ws.Cells.Item($row, 4).Value = $data.D
        $ws.Cells.Item($row, 5).Value = $data.E
        $ws.Cells.Item($row, 6).Value = $data.F
        $ws.Cells.Item($row, 7).Value = $data.G
        $ws.Cells.Item($row, 8).Value = $data.H
        $ws.Cells.Item($row, 9).Value = $data.I
    }

    # Remove the now-duplicated trailing row entirely (was row 5), shrinking
    # the sheet's used range down to row 4.
    $ws.Rows.Item(5).Delete()
}
